$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header + first columns (keep original three fields) ---
$ws.Range("A1").Value = "id"
$ws.Range("C1").Value = "purpose"
$ws.Range("B1").Value = "amount"
$ws.Range("C2").Value = "Ausfallbeitrag"

# --- New "recipient" column (replaces old "Clubnummer") ---
$ws.Range("D1").Value = "recipient"
$ws.Range("D2").Value = "ABC"

# --- New "iban" / "bic" columns ---
$ws.Range("E1").Value = "iban"
$ws.Range("F1").Value = "bic"
$ws.Range("E2").Value = "DE1234"
$ws.Range("F2").Value = "MLPDEBXX"

# --- New "mandate" column ---
$ws.Range("G1").Value = "mandate"
$ws.Range("G2").Value = "ABCDEF"

# --- New header cells (E1:G1) should be bold, matching the rest of row 1 ---
$ws.Range("E1:G1").Font.Bold = $true

# --- Numeric id / amount values ---
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 1

# --- Column widths: drop old bestFit widths on B/D, size C to fit "purpose" ---
$ws.Columns("B").ColumnWidth = 9.67
$ws.Columns("D").ColumnWidth = 9.67
$ws.Columns("C").ColumnWidth = 20.59

# --- Selection moves to A2 ---
$ws.Range("A2").Select()
